# Add two new debt rows (21 and 22, i.e. sheet rows 22 and 23) to the
# "DANH SACH NO" sheet, add a quote-prefixed blank marker in B24, update the
# "THONG KE NAP " log sheet with the matching two new payment entries, and
# leave the selections on both sheets where the author left them.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # DANH SACH NO
$ws2 = $wb.Worksheets.Item(2)   # THONG KE NAP

# ---- Sheet 1: new debt row 22 (STT 21) - Le Ngoc Nhu Y / Mua the Viettel ----
$ws1.Range("B22").Value = "Lê Ngọc Như Ý"
$ws1.Range("C22").Value = "Mua thẻ Viettel"
$ws1.Range("D22").Value = 20000
$ws1.Range("E22").Value = 0
$ws1.Range("G22").Value = 0
$ws1.Range("H22").Value = 0
$ws1.Range("J22").Value = 46024
$ws1.Range("K22").Value = 46030
$ws1.Range("M22").Value = "Chưa trả đủ"

# ---- Sheet 1: new debt row 23 (STT 22) - Nguyen Huu Nhan / Mua the Zing ----
$ws1.Range("B23").Value = "Nguyễn Hữu Nhân"
$ws1.Range("C23").Value = "Mua thẻ Zing"
$ws1.Range("D23").Value = 20000
$ws1.Range("E23").Value = 0
$ws1.Range("G23").Value = 0
$ws1.Range("H23").Value = 0
$ws1.Range("J23").Value = 46024
$ws1.Range("K23").Value = 46030
$ws1.Range("M23").Value = "Chưa trả đủ"

# Fill both formula columns as one shared-formula range, matching the
# pattern already used for the other rows in the table.
$ws1.Range("F22:F23").Formula = "=(D22+I22)-E22"
$ws1.Range("I22:I23").Formula = "=D22*H22"

# ---- Sheet 1: row 24 gets a lone quote-prefixed empty text marker in B24 ----
$ws1.Range("B24").Value = "'"

# ---- Sheet 2: log the two new payments in the next two empty rows ----
$ws2.Range("A137").Value = 46024
$ws2.Range("B137").Value = "Lê Ngọc Như Ý"
$ws2.Range("C137").Value = 20000
$ws2.Range("D137").Value = "Mua thẻ Viettel"

$ws2.Range("A138").Value = 46024
$ws2.Range("B138").Value = "Nguyễn Hữu Nhân"
$ws2.Range("C138").Value = 20000
$ws2.Range("D138").Value = "Mua thẻ Zing"

# ---- Restore on-screen selections the way the author left them ----
$ws2.Activate()
$ws2.Range("A139").Select()

$ws1.Activate()
$ws1.Range("F26").Select()
